$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Runner")

# Update the OperationType values in column B (rows 2-5) to upper-case
# values, matching the refreshed set of allowed operations.
$ws.Range("B2").Value = "INSERT"
$ws.Range("B3").Value = "UPDATE"
$ws.Range("B4").Value = "INSERT"
$ws.Range("B5").Value = "DELETE"

# Add data validation drop-downs for the "Required" (A2:A5) and
# "OperationType" (B2:B5) columns.
$rngA = $ws.Range("A2:A5")
$rngA.Validation.Delete()
$rngA.Validation.Add(3, 1, 1, '"Yes,No"')
$rngA.Validation.IgnoreBlank = $true
$rngA.Validation.InCellDropdown = $true
$rngA.Validation.ShowInput = $true
$rngA.Validation.ShowError = $true

$rngB = $ws.Range("B2:B5")
$rngB.Validation.Delete()
$rngB.Validation.Add(3, 1, 1, '"INSERT,UPDATE,DELETE"')
$rngB.Validation.IgnoreBlank = $true
$rngB.Validation.InCellDropdown = $true
$rngB.Validation.ShowInput = $true
$rngB.Validation.ShowError = $true

# Update the active selection shown when the file is reopened.
$ws.Range("D11").Select()
